$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 49 ("Stack 1D spectra" / "SK"), shifting existing
# rows 49+ down by one.
$ws.Rows.Item(49).Insert()

# Match the row height formatting used by the surrounding shortcut rows.
$ws.Rows.Item(49).RowHeight = 17

# Populate the new shortcut entry.
$ws.Cells.Item(49, 1).Value = "Stack 1D spectra"
$ws.Cells.Item(49, 2).Value = "SK"

# Extend the print area to cover the newly added row.
$ws.PageSetup.PrintArea = '$A$1:$C$127'

# Update the active selection to reflect the shifted position.
$ws.Range("C50").Select()
